# Apply calibration refresh to the WTI riskDriverType PnL factor calibrations workbook.
# Updates scalar parameter values and the autocorrelation-list strings on each
# of the AR / SETAR / GARCH / TARCH / AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.006825207394459255
$ws.Range("B3").Value = 0.7740120446426928
$ws.Range("B4").Value = 0.1066737348607805
$ws.Range("B5").Value = "[1.0, 0.03591633023658923, -0.03202675534859124, -0.10176320494358693, -0.11019256216814352, 0.07496515225145116, 0.018140790266992772, 0.030603399237779473, 0.011666100456573084, -0.018714686383796822, -0.04237695261679498, -0.03351079806106292, -0.008584719366121306, 0.05106730503852144, 0.04853627039262523, 0.03580689289578694, -0.00538901046711703, -0.03770368103526971, -0.03574853699430846, -0.0042334470341165754]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1779548845695258
$ws.Range("B5").Value = 0.5618703997551836
$ws.Range("B6").Value = 0.06552252281878544
$ws.Range("B7").Value = 0.1806170955402093
$ws.Range("B8").Value = 0.5738797443335716
$ws.Range("B9").Value = 0.0623646825431965
$ws.Range("B10").Value = "[0.9999999999999998, 0.08479294654618728, 0.01881629546687745, 0.005265668258526909, -0.004825981679042169, 0.001743930812811218, 0.006134634894567925, -0.02707074731679678, 0.00013728442882836723, -0.01613409079524153, 0.008705904118283434, 0.012609553881461523, -0.003237097404894301, 0.025518956925377815, 0.020936588226391134, 0.0002019237245177806, -0.02130525918623366, -0.01751037344892199, -2.416205482189091e-05, -0.018003619111308346]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.001210009391988734
$ws.Range("B3").Value = 0.112714886432367
$ws.Range("B4").Value = 0.05963541914630002
$ws.Range("B5").Value = [double]"5.989069439454051E-14"
$ws.Range("B6").Value = "[1.0, -0.002267298700155771, -0.028888846210529506, -0.10057383077991917, -0.12177658988596707, 0.09792843714912494, 0.0353070001456782, 0.0399430116107057, 0.015415262220788059, -0.014373956866661741, -0.054595642143399514, -0.02849636875805811, -0.004467030271756254, 0.05169600877879758, 0.038008868427548345, 0.02073322594082356, 0.0009768723167918648, -0.037133387736442036, -0.049227827358005564, -0.010584493127810877]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.00105386075390126
$ws.Range("B3").Value = 0.112623095055708
$ws.Range("B4").Value = 0.02104708133315268
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[0.9999999999999998, -0.0014456543209884585, -0.029892409040372762, -0.10059671877443724, -0.12308192803867783, 0.09797014636914, 0.03582886573049426, 0.03995540255943347, 0.016239788121979654, -0.014453502975154778, -0.05458894249800588, -0.02887284981659763, -0.004567458867819179, 0.05195265317514172, 0.03838938986279587, 0.020493590686416603, 0.0011455914268146234, -0.037926226471186764, -0.04939571790543775, -0.01140523618028735]"
$ws.Range("B7").Value = 0.07857961461559494

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.007186473657713147
$ws.Range("B3").Value = 0.1019441547454095
$ws.Range("B4").Value = 0.00235527575045286
$ws.Range("B5").Value = [double]"3.911594282837372E-09"
$ws.Range("B6").Value = "[1.0, -0.00512539476321441, -0.03429385765695796, -0.10343142142807313, -0.11426245070439539, 0.07830368694206877, 0.015967829928175505, 0.030162640058722677, 0.01153932504129344, -0.01731615938468597, -0.0395013554744115, -0.03134903721579623, -0.008068273968623705, 0.05062079974583103, 0.04581526262333217, 0.03355951469564344, -0.005924246319156136, -0.03715577842559043, -0.03531137514010312, -0.005208636952932974]"
$ws.Range("B7").Value = 0.07904812661771926
$ws.Range("B9").Value = 0.7748920970977445
